$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Développeur" -> "Développement"
$ws.Range("A10").Value = "Développement"

# Row 9: "Chef de projet" -> "Pôle Validation", quantity 3 -> 1
$ws.Range("A9").Value = "Pôle Validation"
$ws.Range("B9").Value = 1

# Move the active selection to B10 (cosmetic cursor position change)
$ws.Range("B10").Select() | Out-Null
